# Regenerate the "K" (strikeouts) column (column G) in the save_data sheet
# so it holds actual strikeout totals per game instead of the old "Strike#"
# (total strikes thrown) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values for rows 2-35 (games 0-33), replacing the old
# "Strike#" pitch counts that used to live in column G.
$kValues = @{
    2  = 1
    3  = 5
    4  = 9
    5  = 3
    6  = 8
    7  = 13
    8  = 10
    9  = 4
    10 = 4
    11 = 8
    12 = 1
    13 = 4
    14 = 4
    15 = 6
    16 = 5
    17 = 7
    18 = 10
    19 = 8
    20 = 5
    21 = 9
    22 = 9
    23 = 5
    24 = 10
    25 = 7
    26 = 6
    27 = 9
    28 = 5
    29 = 6
    30 = 6
    31 = 1
    32 = 0
    33 = 7
    34 = 3
    35 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
